# Fix the 2050 column-header label (was mistakenly left as a leftover numeric
# value "702.2551432549269" instead of the year label) and drop the bottom
# "Total" row from every scenario table.

$wb = $excel.ActiveWorkbook

# Sheets whose 5th column (E) header should simply read "2050"
$simpleYearSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($name in $simpleYearSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("E1").Value = "'2050"
}

# "Potencia Incremental - SIN(MW)" uses period ranges (2015-2030, 2031-2040, ...)
$wsPeriodo = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$wsPeriodo.Range("E1").Value = "'2041-2050"

# Remove the trailing "Total" row (row 13) from the four 13-row tables
$tablesWithTotal = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Potencia Incremental - SIN(MW)"
)

foreach ($name in $tablesWithTotal) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Rows.Item(13).Delete()
}

# Remove the trailing "Total" row (row 4) from "Custo Total (bilhões de R$)"
$wsCusto = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$wsCusto.Rows.Item(4).Delete()
